$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 3
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = 3
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("F21").Value = -3
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = 4
$ws.Range("F25").Value = 4
$ws.Range("F26").Value = -1
$ws.Range("F28").Value = 4
$ws.Range("F34").Value = -4
$ws.Range("F39").Value = 5
$ws.Range("F47").Value = -3
$ws.Range("F50").Value = -2
$ws.Range("F51").Value = -3
$ws.Range("F53").Value = -3
$ws.Range("F55").Value = -4
$ws.Range("F63").Value = -2
$ws.Range("F70").Value = -3
$ws.Range("F72").Value = -2
$ws.Range("F80").Value = -2
